$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-8: in-place odds value updates ---
# Row 2
$ws.Range("Q2").Value = 3.6
$ws.Range("R2").Value = 1.29

# Row 3
$ws.Range("Q3").Value = 2.06
$ws.Range("R3").Value = 1.84

# Row 4
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1.75
$ws.Range("Z4").Value = 12
$ws.Range("AC4").Value = 9
$ws.Range("AK4").Value = 41
$ws.Range("AX4").Value = 29

# Row 5
$ws.Range("G5").Value = 1.7
$ws.Range("I5").Value = 5
$ws.Range("K5").Value = 2.1
$ws.Range("Z5").Value = 13
$ws.Range("AH5").Value = 23
$ws.Range("AU5").Value = 8.5

# Row 6
$ws.Range("G6").Value = 2.6
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 2.63
$ws.Range("AG6").Value = 9
$ws.Range("AN6").Value = 4.75
$ws.Range("AQ6").Value = 51

# Row 8
$ws.Range("G8").Value = 1.91
$ws.Range("I8").Value = 3.9
$ws.Range("J8").Value = 2.6
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("O8").Value = 1.3
$ws.Range("P8").Value = 3.4
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.85
$ws.Range("S8").Value = 1.4
$ws.Range("T8").Value = 2.75
$ws.Range("AC8").Value = 9.5
$ws.Range("AE8").Value = 15
$ws.Range("AG8").Value = 11
$ws.Range("AO8").Value = 10
$ws.Range("AQ8").Value = 34
$ws.Range("AT8").Value = 2.75

# --- Rows 9, 11, 12, 13: full row replacement (new matches / shifted rows) ---
# Row 9
$ws.Range("A9").Value = "jRlzzGEi"
$ws.Range("B9").Value = "09/11/2024"
$ws.Range("C9").Value = "20:05"
$ws.Range("D9").Value = "MEXICO - LIGA MX"
$ws.Range("E9").Value = "Guadalajara Chivas"
$ws.Range("F9").Value = "Atl. San Luis"
$ws.Range("G9").Value = 1.67
$ws.Range("H9").Value = 4.5
$ws.Range("I9").Value = 4.33
$ws.Range("J9").Value = 2.2
$ws.Range("K9").Value = 2.4
$ws.Range("L9").Value = 4.75
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 17
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 4.5
$ws.Range("Q9").Value = 1.6
$ws.Range("R9").Value = 2.3
$ws.Range("S9").Value = 1.3
$ws.Range("T9").Value = 3.4
$ws.Range("U9").Value = 1.67
$ws.Range("V9").Value = 2.1
$ws.Range("W9").Value = 9
$ws.Range("X9").Value = 9
$ws.Range("Y9").Value = 8.5
$ws.Range("Z9").Value = 13
$ws.Range("AA9").Value = 12
$ws.Range("AB9").Value = 21
$ws.Range("AC9").Value = 17
$ws.Range("AD9").Value = 8.5
$ws.Range("AE9").Value = 15
$ws.Range("AF9").Value = 41
$ws.Range("AG9").Value = 15
$ws.Range("AH9").Value = 23
$ws.Range("AI9").Value = 13
$ws.Range("AJ9").Value = 51
$ws.Range("AK9").Value = 34
$ws.Range("AL9").Value = 34
$ws.Range("AM9").Value = 151
$ws.Range("AN9").Value = 3.75
$ws.Range("AO9").Value = 8.5
$ws.Range("AP9").Value = 17
$ws.Range("AQ9").Value = 26
$ws.Range("AR9").Value = 41
$ws.Range("AS9").Value = 101
$ws.Range("AT9").Value = 3.4
$ws.Range("AU9").Value = 7.5
$ws.Range("AV9").Value = 41
$ws.Range("AW9").Value = 6.5
$ws.Range("AX9").Value = 23
$ws.Range("AY9").Value = 26
$ws.Range("AZ9").Value = 81
$ws.Range("BA9").Value = 81
$ws.Range("BB9").Value = 151
$ws.Range("BC9").Value = 451
$ws.Range("BD9").Value = 151

# Row 11
$ws.Range("A11").Value = "jeQ89k4T"
$ws.Range("B11").Value = "09/11/2024"
$ws.Range("C11").Value = "19:00"
$ws.Range("D11").Value = "URUGUAY - PRIMERA DIVISION"
$ws.Range("E11").Value = "Penarol"
$ws.Range("F11").Value = "Liverpool M."
$ws.Range("G11").Value = 1.4
$ws.Range("H11").Value = 4.5
$ws.Range("I11").Value = 7.5
$ws.Range("J11").Value = 1.95
$ws.Range("K11").Value = 2.3
$ws.Range("L11").Value = 7.5
$ws.Range("M11").Value = 1.06
$ws.Range("N11").Value = 10
$ws.Range("O11").Value = 1.29
$ws.Range("P11").Value = 3.5
$ws.Range("Q11").Value = 1.98
$ws.Range("R11").Value = 1.88
$ws.Range("S11").Value = 1.4
$ws.Range("T11").Value = 2.75
$ws.Range("U11").Value = 2.25
$ws.Range("V11").Value = 1.57
$ws.Range("W11").Value = 6
$ws.Range("X11").Value = 6
$ws.Range("Y11").Value = 9
$ws.Range("Z11").Value = 8.5
$ws.Range("AA11").Value = 13
$ws.Range("AB11").Value = 34
$ws.Range("AC11").Value = 10
$ws.Range("AD11").Value = 9
$ws.Range("AE11").Value = 23
$ws.Range("AF11").Value = 81
$ws.Range("AG11").Value = 15
$ws.Range("AH11").Value = 41
$ws.Range("AI11").Value = 23
$ws.Range("AJ11").Value = 81
$ws.Range("AK11").Value = 51
$ws.Range("AL11").Value = 51
$ws.Range("AM11").Value = 1250
$ws.Range("AN11").Value = 3.2
$ws.Range("AO11").Value = 7
$ws.Range("AP11").Value = 21
$ws.Range("AQ11").Value = 21
$ws.Range("AR11").Value = 51
$ws.Range("AS11").Value = 151
$ws.Range("AT11").Value = 2.75
$ws.Range("AU11").Value = 10
$ws.Range("AV11").Value = 67
$ws.Range("AW11").Value = 8.5
$ws.Range("AX11").Value = 41
$ws.Range("AY11").Value = 41
$ws.Range("AZ11").Value = 201
$ws.Range("BA11").Value = 201
$ws.Range("BB11").Value = 301
$ws.Range("BC11").Value = 51
$ws.Range("BD11").Value = 51

# Row 12
$ws.Range("A12").Value = "prLbSHS6"
$ws.Range("B12").Value = "09/11/2024"
$ws.Range("C12").Value = "18:00"
$ws.Range("D12").Value = "USA - MLS"
$ws.Range("E12").Value = "FC Cincinnati"
$ws.Range("F12").Value = "New York City"
$ws.Range("G12").Value = 1.85
$ws.Range("H12").Value = 3.7
$ws.Range("I12").Value = 3.9
$ws.Range("J12").Value = 2.4
$ws.Range("K12").Value = 2.3
$ws.Range("L12").Value = 4.33
$ws.Range("M12").Value = 1.04
$ws.Range("N12").Value = 13
$ws.Range("O12").Value = 1.22
$ws.Range("P12").Value = 4
$ws.Range("Q12").Value = 1.7
$ws.Range("R12").Value = 2.1
$ws.Range("S12").Value = 1.33
$ws.Range("T12").Value = 3.25
$ws.Range("U12").Value = 1.62
$ws.Range("V12").Value = 2.2
$ws.Range("W12").Value = 9
$ws.Range("X12").Value = 10
$ws.Range("Y12").Value = 8.5
$ws.Range("Z12").Value = 17
$ws.Range("AA12").Value = 13
$ws.Range("AB12").Value = 21
$ws.Range("AC12").Value = 13
$ws.Range("AD12").Value = 7
$ws.Range("AE12").Value = 13
$ws.Range("AF12").Value = 41
$ws.Range("AG12").Value = 13
$ws.Range("AH12").Value = 21
$ws.Range("AI12").Value = 13
$ws.Range("AJ12").Value = 41
$ws.Range("AK12").Value = 29
$ws.Range("AL12").Value = 34
$ws.Range("AM12").Value = 151
$ws.Range("AN12").Value = 4
$ws.Range("AO12").Value = 9.5
$ws.Range("AP12").Value = 17
$ws.Range("AQ12").Value = 29
$ws.Range("AR12").Value = 41
$ws.Range("AS12").Value = 101
$ws.Range("AT12").Value = 3.25
$ws.Range("AU12").Value = 7.5
$ws.Range("AV12").Value = 41
$ws.Range("AW12").Value = 6
$ws.Range("AX12").Value = 21
$ws.Range("AY12").Value = 23
$ws.Range("AZ12").Value = 67
$ws.Range("BA12").Value = 81
$ws.Range("BB12").Value = 151
$ws.Range("BC12").Value = 451
$ws.Range("BD12").Value = 151

# Row 13
$ws.Range("A13").Value = "KGLtPqC6"
$ws.Range("B13").Value = "09/11/2024"
$ws.Range("C13").Value = "20:00"
$ws.Range("D13").Value = "USA - MLS"
$ws.Range("E13").Value = "Orlando City"
$ws.Range("F13").Value = "Charlotte"
$ws.Range("G13").Value = 1.8
$ws.Range("H13").Value = 3.75
$ws.Range("I13").Value = 4.2
$ws.Range("J13").Value = 2.38
$ws.Range("K13").Value = 2.25
$ws.Range("L13").Value = 4.5
$ws.Range("M13").Value = 1.04
$ws.Range("N13").Value = 13
$ws.Range("O13").Value = 1.22
$ws.Range("P13").Value = 4
$ws.Range("Q13").Value = 1.8
$ws.Range("R13").Value = 2
$ws.Range("S13").Value = 1.36
$ws.Range("T13").Value = 3
$ws.Range("U13").Value = 1.75
$ws.Range("V13").Value = 2
$ws.Range("W13").Value = 8
$ws.Range("X13").Value = 9
$ws.Range("Y13").Value = 8.5
$ws.Range("Z13").Value = 15
$ws.Range("AA13").Value = 15
$ws.Range("AB13").Value = 23
$ws.Range("AC13").Value = 12
$ws.Range("AD13").Value = 7
$ws.Range("AE13").Value = 15
$ws.Range("AF13").Value = 51
$ws.Range("AG13").Value = 13
$ws.Range("AH13").Value = 21
$ws.Range("AI13").Value = 13
$ws.Range("AJ13").Value = 41
$ws.Range("AK13").Value = 34
$ws.Range("AL13").Value = 34
$ws.Range("AM13").Value = 201
$ws.Range("AN13").Value = 3.75
$ws.Range("AO13").Value = 9.5
$ws.Range("AP13").Value = 19
$ws.Range("AQ13").Value = 29
$ws.Range("AR13").Value = 51
$ws.Range("AS13").Value = 126
$ws.Range("AT13").Value = 3
$ws.Range("AU13").Value = 8
$ws.Range("AV13").Value = 51
$ws.Range("AW13").Value = 6
$ws.Range("AX13").Value = 23
$ws.Range("AY13").Value = 29
$ws.Range("AZ13").Value = 81
$ws.Range("BA13").Value = 81
$ws.Range("BB13").Value = 201
$ws.Range("BC13").Value = 501
$ws.Range("BD13").Value = 151

